# Auto-applies the scheduled-runner price/profit refresh described in the commit diff.
# For each (sheet, row) touched, update the changed numeric cells; where a cell was
# removed outright in the diff (profit columns with no sellable HQ/NQ data), clear it
# instead of writing a stale 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 302.6216
$ws.Range("I15").Value = 302.6216
$ws.Range("K15").Value = 907.8648000000001
$ws.Range("M15").Value = -738.8648000000001

$ws.Range("H137").Value = 2171.5833
$ws.Range("I137").Value = 1632.579
$ws.Range("K137").Value = 4897.737
$ws.Range("M137").Value = -2347.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1575
$ws.Range("I61").Value = 1575
$ws.Range("K61").Value = 1575
$ws.Range("M61").Value = -1363

$ws.Range("H74").Value = 733
$ws.Range("I74").Value = 392.5
$ws.Range("K74").Value = 392.5
$ws.Range("M74").Value = 481.5

$ws.Range("H77").Value = 733
$ws.Range("I77").Value = 392.5
$ws.Range("K77").Value = 1962.5
$ws.Range("M77").Value = 2405.5

$ws.Range("H102").Value = 2060.5
$ws.Range("I102").Value = 1247.6666
$ws.Range("K102").Value = 1247.6666
$ws.Range("M102").Value = 374.3334

$ws.Range("H136").Value = 1575
$ws.Range("I136").Value = 1575
$ws.Range("K136").Value = 4725
$ws.Range("M136").Value = -2175

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 419.16666
$ws.Range("J80").Value = 355.57144
$ws.Range("L80").Value = 355.57144
$ws.Range("N80").Value = -2351.57144

$ws.Range("H82").Value = 5144
$ws.Range("I82").Value = 5144
$ws.Range("K82").Value = 5144
$ws.Range("M82").Value = -4761

$ws.Range("H83").Value = 419.16666
$ws.Range("J83").Value = 355.57144
$ws.Range("L83").Value = 1777.8572
$ws.Range("N83").Value = -11761.8572

$ws.Range("H85").Value = 5144
$ws.Range("I85").Value = 5144
$ws.Range("K85").Value = 5144
$ws.Range("M85").Value = -3818

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 5252
$ws.Range("J134").Value = 3200
$ws.Range("L134").Value = 9600
$ws.Range("N134").Value = -14670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1811.762
$ws.Range("I31").Value = 1863.8889
$ws.Range("K31").Value = 1863.8889
$ws.Range("M31").Value = -1568.8889

$ws.Range("H34").Value = 1811.762
$ws.Range("I34").Value = 1863.8889
$ws.Range("K34").Value = 1863.8889
$ws.Range("M34").Value = -1661.8889

$ws.Range("H58").Value = 1350.5555
$ws.Range("I58").Value = 1028.8667
$ws.Range("K58").Value = 1028.8667
$ws.Range("M58").Value = -825.8667

$ws.Range("H132").Value = 1865.909
$ws.Range("I132").Value = 1988.9656
$ws.Range("K132").Value = 5966.8968
$ws.Range("M132").Value = -3436.8968

$ws.Range("H136").Value = 1350.5555
$ws.Range("I136").Value = 1028.8667
$ws.Range("K136").Value = 3086.6001
$ws.Range("M136").Value = -536.6001000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1109.9
$ws.Range("J131").Value = 1111.0555
$ws.Range("L131").Value = 3333.1665
$ws.Range("N131").Value = -13413.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 680
$ws.Range("I2").Value = 2026.6
$ws.Range("J2").Value = 283.94116
$ws.Range("K2").Value = 2026.6
$ws.Range("L2").Value = 283.94116
$ws.Range("M2").Value = -1913.6
$ws.Range("N2").Value = -509.94116

$ws.Range("H101").Value = 34528.5
$ws.Range("J101").Value = 34528.5
$ws.Range("L101").Value = 34528.5
$ws.Range("N101").Value = -41018.5

$ws.Range("H102").Value = 1397
$ws.Range("I102").Value = 1150
$ws.Range("K102").Value = 1150
$ws.Range("M102").Value = 472

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1769.1482
$ws.Range("J46").Value = 2218.7273
$ws.Range("L46").Value = 2218.7273
$ws.Range("N46").Value = -2594.7273

$ws.Range("H55").Value = 449.92307
$ws.Range("J55").Value = 609.7778
$ws.Range("L55").Value = 609.7778
$ws.Range("N55").Value = -955.7778

$ws.Range("H93").Value = 1850.5
$ws.Range("I93").Value = 1850.5
$ws.Range("K93").Value = 1850.5
$ws.Range("M93").Value = -602.5

$ws.Range("H122").Value = 5605.5
$ws.Range("I122").Value = 5581.5
$ws.Range("J122").Value = 5641.5
$ws.Range("K122").Value = 16744.5
$ws.Range("L122").Value = 16924.5
$ws.Range("M122").Value = -14294.5
$ws.Range("N122").Value = -21824.5

$ws.Range("H132").Value = 4940.3
$ws.Range("I132").Value = 5135.706
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 15407.118
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -12877.118
$ws.Range("N132").Value = -16559

$ws.Range("H136").Value = 4717.727
$ws.Range("I136").Value = 3985
$ws.Range("K136").Value = 11955
$ws.Range("M136").Value = -9405

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5748.75
$ws.Range("I62").Value = 5999
$ws.Range("J62").Value = 5498.5
$ws.Range("K62").Value = 5999
$ws.Range("L62").Value = 5498.5
$ws.Range("M62").Value = -5375
$ws.Range("N62").Value = -6746.5

$ws.Range("H65").Value = 5748.75
$ws.Range("I65").Value = 5999
$ws.Range("J65").Value = 5498.5
$ws.Range("K65").Value = 29995
$ws.Range("L65").Value = 27492.5
$ws.Range("M65").Value = -26875
$ws.Range("N65").Value = -33732.5

$ws.Range("H81").Value = 6812.125
$ws.Range("I81").Value = 4166.6665
$ws.Range("K81").Value = 8333.333000000001
$ws.Range("M81").Value = -7272.333000000001

$ws.Range("H84").Value = 6812.125
$ws.Range("I84").Value = 4166.6665
$ws.Range("K84").Value = 41666.665
$ws.Range("M84").Value = -36362.665

$ws.Range("H96").Value = 3185.5715
$ws.Range("I96").Value = 2899.75
$ws.Range("K96").Value = 2899.75
$ws.Range("M96").Value = -1526.75

$ws.Range("H107").Value = 301
$ws.Range("I107").Value = 331.14285
$ws.Range("K107").Value = 993.4285500000001
$ws.Range("M107").Value = 926.5714499999999

$ws.Range("H132").Value = 640.7059
$ws.Range("I132").Value = 640.7059
$ws.Range("K132").Value = 1922.1177
$ws.Range("M132").Value = 607.8822999999998

$ws.Range("H136").Value = 2538.75
$ws.Range("I136").Value = 1758.7142
$ws.Range("K136").Value = 5276.142599999999
$ws.Range("M136").Value = -2726.142599999999
